$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1423941644592333
$ws.Range("D2").Value = 0.05807279407734711
$ws.Range("E2").Value = 0.4087874611668241
$ws.Range("F2").Value = 1.79638570722858
$ws.Range("G2").Value = 0.0024656531813163
$ws.Range("K2").Value = 1.694994303719
$ws.Range("N2").Value = 1.862483254517457

$ws.Range("B3").Value = 0.1329530972492705
$ws.Range("D3").Value = 0.05817162727246128
$ws.Range("E3").Value = 0.3562932915316708
$ws.Range("F3").Value = 1.735381256389203
$ws.Range("G3").Value = 0.002471821482353404
$ws.Range("K3").Value = 1.518540894162186
$ws.Range("N3").Value = 1.870822257611934

$ws.Range("B4").Value = 0.1272303690788164
$ws.Range("D4").Value = 0.05826370576566831
$ws.Range("E4").Value = 0.3242131329484721
$ws.Range("F4").Value = 1.699164556402764
$ws.Range("G4").Value = 0.002475799153064283
$ws.Range("K4").Value = 1.411087437083324
$ws.Range("N4").Value = 1.876562783581662

$ws.Range("B5").Value = 0.1249170176737806
$ws.Range("D5").Value = 0.05830903056088843
$ws.Range("E5").Value = 0.3111749585332717
$ws.Range("F5").Value = 1.684714150371974
$ws.Range("G5").Value = 0.002477468136768102
$ws.Range("K5").Value = 1.367516525671363
$ws.Range("N5").Value = 1.879057246869735

$ws.Range("B6").Value = 0.1245340193066511
$ws.Range("D6").Value = 0.05831702526318949
$ws.Range("E6").Value = 0.3090119818158712
$ws.Range("F6").Value = 1.682333178880455
$ws.Range("G6").Value = 0.002477748178305293
$ws.Range("K6").Value = 1.360294543318332
$ws.Range("N6").Value = 1.879480797056914

$ws.Range("B7").Value = 0.1271990945780601
$ws.Range("D7").Value = 0.0582642855612967
$ws.Range("E7").Value = 0.3240371592440425
$ws.Range("F7").Value = 1.698968430134329
$ws.Range("G7").Value = 0.002475821466867762
$ws.Range("K7").Value = 1.410498953382671
$ws.Range("N7").Value = 1.876595797675492

$ws.Range("B8").Value = 0.1391235364462773
$ws.Range("D8").Value = 0.05810029969813257
$ws.Range("E8").Value = 0.3906539952307782
$ws.Range("F8").Value = 1.775091696737846
$ws.Range("G8").Value = 0.002467740631128563
$ws.Range("K8").Value = 1.633964544863545
$ws.Range("N8").Value = 1.865229282749311

$ws.Range("B9").Value = 0.1630939545643884
$ws.Range("D9").Value = 0.05803208871130749
$ws.Range("E9").Value = 0.5226494177196344
$ws.Range("F9").Value = 1.934387811767778
$ws.Range("G9").Value = 0.00245339517264902
$ws.Range("K9").Value = 2.079556702365721
$ws.Range("N9").Value = 1.847900062226273

$ws.Range("B10").Value = 0.1810627156776974
$ws.Range("D10").Value = 0.05814248345022577
$ws.Range("E10").Value = 0.6206855936165994
$ws.Range("F10").Value = 2.057792685750201
$ws.Range("G10").Value = 0.0024437580025164
$ws.Range("K10").Value = 2.411908645233268
$ws.Range("N10").Value = 1.838246175388875

$ws.Range("B11").Value = 0.1893150198794302
$ws.Range("D11").Value = 0.05822890746777176
$ws.Range("E11").Value = 0.6655660802868368
$ws.Range("F11").Value = 2.11537466235697
$ws.Range("G11").Value = 0.002439567040554846
$ws.Range("K11").Value = 2.564292154730026
$ws.Range("N11").Value = 1.834534601588018

$ws.Range("B12").Value = 0.1924511733605385
$ws.Range("D12").Value = 0.05826695333356469
$ws.Range("E12").Value = 0.6826060028226806
$ws.Range("F12").Value = 2.137391787888248
$ws.Range("G12").Value = 0.002438007579795179
$ws.Range("K12").Value = 2.622175948443328
$ws.Range("N12").Value = 1.833227910256056

$ws.Range("B13").Value = 0.1917752501320962
$ws.Range("D13").Value = 0.05825852108804241
$ws.Range("E13").Value = 0.6789341013130468
$ws.Range("F13").Value = 2.132640503107922
$ws.Range("G13").Value = 0.002438342214528462
$ws.Range("K13").Value = 2.609701532520205
$ws.Range("N13").Value = 1.833504919059351

$ws.Range("B14").Value = 0.1895728093075491
$ws.Range("D14").Value = 0.05823193026785844
$ws.Range("E14").Value = 0.6669670472663967
$ws.Range("F14").Value = 2.117181748778847
$ws.Range("G14").Value = 0.002439438191544833
$ws.Range("K14").Value = 2.569050648138614
$ws.Range("N14").Value = 1.834425113315717

$ws.Range("B15").Value = 0.1882252054682283
$ws.Range("D15").Value = 0.0582163387019321
$ws.Range("E15").Value = 0.6596428142970581
$ws.Range("F15").Value = 2.107740565696048
$ws.Range("G15").Value = 0.002440113092942624
$ws.Range("K15").Value = 2.544174409094182
$ws.Range("N15").Value = 1.835001656078362

$ws.Range("B16").Value = 0.1805249767861028
$ws.Range("D16").Value = 0.05813757354859916
$ws.Range("E16").Value = 0.6177585822569824
$ws.Range("F16").Value = 2.054059005952809
$ws.Range("G16").Value = 0.00244403575950625
$ws.Range("K16").Value = 2.40197471836899
$ws.Range("N16").Value = 1.838502505632349

$ws.Range("B17").Value = 0.1758211252122379
$ws.Range("D17").Value = 0.05809860598525773
$ws.Range("E17").Value = 0.5921392891548294
$ws.Range("F17").Value = 2.021500151810017
$ws.Range("G17").Value = 0.002446491492306622
$ws.Range("K17").Value = 2.3150512434093
$ws.Range("N17").Value = 1.840825050703998

$ws.Range("B18").Value = 0.1731229607540428
$ws.Range("D18").Value = 0.05807959374582339
$ws.Range("E18").Value = 0.5774301793740904
$ws.Range("F18").Value = 2.00290885484813
$ws.Range("G18").Value = 0.002447922145749062
$ws.Range("K18").Value = 2.265167166055221
$ws.Range("N18").Value = 1.842224862960578

$ws.Range("B19").Value = 0.1722106758140569
$ws.Range("D19").Value = 0.05807373723857268
$ws.Range("E19").Value = 0.5724543491298704
$ws.Range("F19").Value = 1.996637335274158
$ws.Range("G19").Value = 0.002448409668827868
$ws.Range("K19").Value = 2.248296282105514
$ws.Range("N19").Value = 1.842709767454352

$ws.Range("B20").Value = 0.1763210962985653
$ws.Range("D20").Value = 0.05810240136357692
$ws.Range("E20").Value = 0.5948637431657176
$ws.Range("F20").Value = 2.024952022186511
$ws.Range("G20").Value = 0.002446228195148174
$ws.Range("K20").Value = 2.324292748737832
$ws.Range("N20").Value = 1.840571186387493

$ws.Range("B21").Value = 0.1902194165059257
$ws.Range("D21").Value = 0.05823959536892431
$ws.Range("E21").Value = 0.6704808175591523
$ws.Range("F21").Value = 2.121716568588141
$ws.Range("G21").Value = 0.002439115529828838
$ws.Range("K21").Value = 2.580985865887214
$ws.Range("N21").Value = 1.834152140351193

$ws.Range("B22").Value = 0.1993679386963407
$ws.Range("D22").Value = 0.05836031687433518
$ws.Range("E22").Value = 0.7201631161930493
$ws.Range("F22").Value = 2.18619605412448
$ws.Range("G22").Value = 0.002434627575358192
$ws.Range("K22").Value = 2.749799393632316
$ws.Range("N22").Value = 1.830533294493634

$ws.Range("B23").Value = 0.1944792579969032
$ws.Range("D23").Value = 0.05829300582215069
$ws.Range("E23").Value = 0.6936214344406579
$ws.Range("F23").Value = 2.151667371610046
$ws.Range("G23").Value = 0.002437008251489494
$ws.Range("K23").Value = 2.65960181620801
$ws.Range("N23").Value = 1.832411662877448

$ws.Range("B24").Value = 0.1760950402325392
$ws.Range("D24").Value = 0.05810067492046755
$ws.Range("E24").Value = 0.5936319562445931
$ws.Range("F24").Value = 2.023391035650263
$ws.Range("G24").Value = 0.002446347173128638
$ws.Range("K24").Value = 2.320114386335433
$ws.Range("N24").Value = 1.840685757462211

$ws.Range("B25").Value = 0.1565465185942259
$ws.Range("D25").Value = 0.05802280047137032
$ws.Range("E25").Value = 0.4867720553613708
$ws.Range("F25").Value = 1.890194602601071
$ws.Range("G25").Value = 0.00245711660854621
$ws.Range("K25").Value = 1.958173624638732
$ws.Range("N25").Value = 1.852051513601197
